# Rename the three inline picture shapes (Pearson Edexcel logo x2 in the
# footers, BTEC logo x1 in the header) so their docPr/cNvPr "name" labels
# swap between "image1"/"image2" (the embedded media files themselves are
# untouched - only the display name metadata on each picture changes).
#
# InlineShape has no settable Name in the Word object model, so each
# picture is temporarily converted to a floating Shape (which does expose
# a settable .Name), renamed, and converted back to an InlineShape so the
# drawing stays <wp:inline> exactly as before.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $newName) {
    $shape = $range.InlineShapes(1)
    $floating = $shape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

# First-page header -> BTEC logo: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers(2).Range "image1.jpg"

# Default footer -> Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers(1).Range "image2.png"

# First-page footer -> Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers(2).Range "image2.png"

Write-Output "Renamed 3 inline picture shapes."
